$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new header columns (AD, AE, AF) mirroring the style of the
# existing header cells (bold, bordered, centered) by copying the format
# from the last existing header cell (AC1).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every data row.
for ($r = 2; $r -le 66; $r++) {
    $ws.Cells.Item($r, 30).Value = 60   # AD
    $ws.Cells.Item($r, 31).Value = 102  # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
